$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.21011
$ws.Range("H2").Value = 0.63033
$ws.Range("M2").Value = 10.25883033333333
$ws.Range("N2").Value = 30.776491
$ws.Range("O2").Value = 0.34684992242997
$ws.Range("P2").Value = 0.34684992242997
$ws.Range("Q2").Value = 2.155482841336667
$ws.Range("R2").Value = 19.39934557203
$ws.Range("S2").Value = 0.34684992242997
$ws.Range("T2").Value = 0.34684992242997

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.21011
$ws.Range("H3").Value = 0.63033
$ws.Range("O3").Value = 0.1682819529322607
$ws.Range("P3").Value = 0.1682819529322608
$ws.Range("Q3").Value = 1.045780432963333
$ws.Range("R3").Value = 9.41202389667
$ws.Range("S3").Value = 0.1682819529322607
$ws.Range("T3").Value = 0.1682819529322608

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.21011
$ws.Range("H4").Value = 0.63033
$ws.Range("M4").Value = 4.130648333333333
$ws.Range("N4").Value = 12.391945
$ws.Range("O4").Value = 0.1396567647041521
$ws.Range("P4").Value = 0.1396567647041521
$ws.Range("Q4").Value = 0.8678905213166667
$ws.Range("R4").Value = 7.811014691849999
$ws.Range("S4").Value = 0.1396567647041521
$ws.Range("T4").Value = 0.1396567647041521

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.21011
$ws.Range("H5").Value = 0.63033
$ws.Range("M5").Value = 2.760918333333333
$ws.Range("N5").Value = 8.282755
$ws.Range("O5").Value = 0.09334634443076846
$ws.Range("P5").Value = 0.09334634443076847
$ws.Range("Q5").Value = 0.5800965510166666
$ws.Range("R5").Value = 5.22086895915
$ws.Range("S5").Value = 0.09334634443076846
$ws.Range("T5").Value = 0.09334634443076847

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.21011
$ws.Range("H6").Value = 0.63033
$ws.Range("M6").Value = 2.602884
$ws.Range("N6").Value = 7.808651999999999
$ws.Range("O6").Value = 0.08800322104565558
$ws.Range("P6").Value = 0.0880032210456556
$ws.Range("Q6").Value = 0.54689195724
$ws.Range("R6").Value = 4.922027615159999
$ws.Range("S6").Value = 0.08800322104565558
$ws.Range("T6").Value = 0.0880032210456556

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.21011
$ws.Range("H7").Value = 0.63033
$ws.Range("M7").Value = 4.846564
$ws.Range("N7").Value = 14.539692
$ws.Range("O7").Value = 0.1638617944571931
$ws.Range("P7").Value = 0.1638617944571932
$ws.Range("Q7").Value = 1.01831156204
$ws.Range("R7").Value = 9.16480405836
$ws.Range("S7").Value = 0.1638617944571931
$ws.Range("T7").Value = 0.1638617944571932
